# Add the new data row (2025-11-20, 92) to the ORA Errors "Online" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 18 (the prior last row) already carries the date number-format style
# used by column A. Copy its formatting down to the new row 19 first so the
# new date cell reuses the same style index instead of Excel creating a
# brand-new (duplicate) style entry.
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)  # xlPasteFormats

# Now fill in the actual values for the new row.
$ws.Range("A19").Value = 45981   # 11/20/2025, stored as a date serial like the rows above it
$ws.Range("B19").Value = 92

# Match the author's final selection/active cell on the newly added row.
$ws.Range("A19:B19").Select()
